$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit message: "Updating excel documents with columns and weights"
# Row 2 (the Erasmus data row) gets its "category" (D) and "title" (E)
# columns filled in, plus a numeric weight under the "study" (J) column.
$ws.Range("D2").Value = "student"
$ws.Range("E2").Value = "Erasmus programs"
$ws.Range("J2").Value = 20

# Reflect the author's final selection/active cell in the saved view state.
[void]$ws.Range("L2").Select()
